$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 5 new columns before column G (shifts existing G:V data to L:AA)
$ws.Range("G1:K1").EntireColumn.Insert()

# 2. New header labels for the inserted "meta" statistic columns
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# 3. Apply the currency number format used by the other monetary columns
$ws.Range("G2:K6").NumberFormat = "R$ #,##0.00"

# 4. Fill in the new "meta" statistic values for every data row
$ws.Range("G2").Value = 11850.23217085951
$ws.Range("H2").Value = 1692.890310122787
$ws.Range("I2").Value = 1741.340023039662
$ws.Range("J2").Value = 13.21358111812448
$ws.Range("K2").Value = 3923.898127062306

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0

$ws.Range("G4").Value = 13609.25823278627
$ws.Range("H4").Value = 972.0898737704478
$ws.Range("I4").Value = 1272.337369789563
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3770.073700501662

$ws.Range("G5").Value = 29904.02695485389
$ws.Range("H5").Value = 1869.001684678368
$ws.Range("I5").Value = 2138.648557943156
$ws.Range("J5").Value = 63.11299939290937
$ws.Range("K5").Value = 8371.809869542562

$ws.Range("G6").Value = 109835.540456445
$ws.Range("H6").Value = 1098.35540456445
$ws.Range("I6").Value = 2287.409798579507
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 21176.91783511972

# 5. Minor recomputation corrections in the shifted (pre-existing) columns
$ws.Range("W2").Value = 2.497617912751116
$ws.Range("L4").Value = 5551.368744100645
$ws.Range("R5").Value = 9.871079671113664
$ws.Range("L6").Value = 26161.59630367917
$ws.Range("Q6").Value = 21.37695663886887
$ws.Range("W6").Value = 30.79468667274806
